$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pages_with_menu")
$ws2 = $wb.Worksheets.Item("mega_menu_mobile")
$new = $wb.Worksheets.Add($ws2)
$new.Name = "sub_menu"

# Column C (index) first
$new.Range("C1").Value = "index"
$cIndex = @(6,5,6,5,6,6,5,6,5,6,5)
for ($i = 0; $i -lt $cIndex.Length; $i++) {
    $new.Cells.Item($i + 2, 3).Value = $cIndex[$i]
}

# Column B (boundary) next
$new.Range("B1").Value = "boundary"
$bBoundary = @(
  "250:500,250:500,350:550,400:600,150:350,200:400",
  "250:450,250:450,250:450,50:150,120:320",
  "250:500,250:500,350:550,400:600,150:350,200:400",
  "250:450,250:450,250:450,50:150,120:320",
  "250:500,250:500,350:550,400:600,150:350,200:400",
  "250:500,250:500,350:550,400:600,150:350,200:400",
  "250:450,250:450,250:450,50:150,120:320",
  "250:500,250:500,350:550,400:600,150:350,200:400",
  "250:450,250:450,250:450,50:150,120:320",
  "250:500,250:500,350:550,400:600,150:350,200:400",
  "250:450,250:450,250:450,50:150,120:320"
)
for ($i = 0; $i -lt $bBoundary.Length; $i++) {
    $new.Cells.Item($i + 2, 2).Value = $bBoundary[$i]
}

# Column A (path) last
$new.Range("A1").Value = "path"
$aPath = @(
  "about-cancer/coping/feelings/relaxation/vitamin-d-supplement-cancer-prevention",
  "espanol/noticias/comunicados-de-prensa/2018/oropharyngeal-hpv-cisplatin",
  "about-cancer/coping/feelings",
  "espanol/cancer/sobrellevar/sentimientos",
  "types/breast/patient/breast-treatment-pdq",
  "news-events/press-releases/2018/oropharyngeal-hpv-cisplatin",
  "espanol/cancer/sobrellevar/sentimientos/hoja-informativa-estres",
  "nano",
  "nano/espanol/cancer/sobrellevar/sentimientos",
  "pediatric-adult-rare-tumor",
  "pediatric-adult-rare-tumor/espanol"
)
for ($i = 0; $i -lt $aPath.Length; $i++) {
    $new.Cells.Item($i + 2, 1).Value = $aPath[$i]
}

# Styles
$new.Range("A1:B1").Font.Bold = $true
$new.Range("C1").Font.Bold = $true
$new.Range("C1").HorizontalAlignment = -4108

$new.Range("C2:C12").Font.Size = 9
$new.Range("C2:C12").Font.Name = "Menlo"

# Column widths
$new.Columns.Item(1).ColumnWidth = 87.1640625
$new.Columns.Item(2).ColumnWidth = 44.1640625

# Selections / active sheet
$ws1.Activate()
$ws1.Range("A1:D13").Select()
$new.Activate()
$new.Range("C14").Select()

Write-Host "done"
